$d = $word.ActiveDocument

# Locate the paragraph that holds the M2Doc field (the one whose instruction
# text ends with ...fromHTMLBodyString('http://www.m2doc.org/tests/')) by
# scanning the paragraphs for the one whose range spans the field's fldChar.
$field = $d.Fields.Item(1)
$fieldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($field.Code.Start -ge $candidate.Range.Start -and $field.Code.Start -lt $candidate.Range.End) {
        $fieldPara = $candidate.Range
        break
    }
}
if ($fieldPara -eq $null) {
    # Fallback: the field-holding paragraph is the 2nd paragraph in this document.
    $fieldPara = $d.Paragraphs.Item(2).Range
}

# Rebuild that paragraph's content, splitting the former single
# "http://www.m2doc.org/tests/" instrText run into "http" + "s" + "://www.m2doc.org/tests/"
# (https instead of http) and moving the "_GoBack" bookmark so it still sits
# right before the "://www.m2doc.org/tests/" run, mirroring where Word would
# leave it after typing the extra "s".
$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">
<w:r><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>
<w:r w:rsidR="00DE6D5A"><w:instrText>m</w:instrText></w:r>
<w:r w:rsidR="002033E1"><w:instrText>:</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText xml:space="preserve"> ('</w:instrText></w:r>
<w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:instrText>&lt;img src="../images/logo_M2Doc.png" alt="" height="54"&gt;</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText>&lt;h2 id="starting-with</w:instrText></w:r>
<w:r w:rsidR="00D62429"><w:instrText>-m2doc"&gt;Starting with ' + self.</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText>n</w:instrText></w:r>
<w:r w:rsidR="00D62429"><w:instrText>a</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText>me + '&lt;/h2&gt;').from</w:instrText></w:r>
<w:r w:rsidR="00342B27"><w:instrText>HTML</w:instrText></w:r>
<w:r w:rsidR="00256E67"><w:instrText>Body</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText>String(</w:instrText></w:r>
<w:r w:rsidR="00D62429"><w:instrText>'</w:instrText></w:r>
<w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:instrText>http</w:instrText></w:r>
<w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:instrText>s</w:instrText></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:instrText>://www.m2doc.org/tests/</w:instrText></w:r>
<w:r w:rsidR="00D62429"><w:instrText>'</w:instrText></w:r>
<w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:instrText>)</w:instrText></w:r>
<w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>
<w:r><w:fldChar w:fldCharType="end"/></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$fieldPara.InsertXML($newParaXml)
